# Generate Report for Handback
# The a94a3e9b-... file has moved from "Ready for handoff" to
# "Handed back: in sync with en-US": update Status/summary cells across the
# Overview sheet and the per-locale (zh-cn / de-de) report sheets, refresh
# the "Latest Handback DateTime" timestamps, and clear the stale
# "Error Detail" message now that the handback is in sync.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for a94a3e9b-c574-4087-9961-b9d0b4140a4e.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row for a94a3e9b-c574-4087-9961-b9d0b4140a4e.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-21 18:52:31"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet: row for a94a3e9b-c574-4087-9961-b9d0b4140a4e.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-21 18:52:38"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
